$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.881.90'
$ws.Range("E2").Value = '  -1.88%  '
$ws.Range("D3").Value = '1.825.25'
$ws.Range("E3").Value = '  -2.20%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.9994'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '240.01'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.39%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.6897'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("E7").Value = '  -0.15%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.07599'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -3.20%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.3016'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -3.96%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '23.37'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07735'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -3.52%  '
$ws.Range("D12").Value = '1.826.85'
$ws.Range("E12").Value = '  -2.30%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '5.036'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -3.19%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '90.01'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -3.62%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.6710'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -4.25%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '6.347'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -1.83%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.000008258'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").Value = '28.865.43'
$ws.Range("E18").Value = '  -2.04%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '242.46'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -4.14%  '
$ws.Range("D20").Value = '2.076.44'
$ws.Range("E20").Value = '  -2.59%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '12.59'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -4.20%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.9998'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '7.383'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -2.97%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '0.9999'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -0.10%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.1469'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -5.61%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '160.75'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '8.703'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -3.48%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '18.14'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.531'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +2.03%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '4.185'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -3.21%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '4.113'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -3.89%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.190'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -1.53%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.05085'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -4.04%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.7508'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -0.11%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.804'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -4.14%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.136'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -2.61%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.680'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.01833'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -2.38%  '
$ws.Range("D39").Value = '1.200.27'
$ws.Range("E39").Value = '  -4.96%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.679'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -2.22%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.9237'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +3.00%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '107.88'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -0.81%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.9996'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.15%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.5162'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").Value = '1.976.53'
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("E46").Value = '  -5.93%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '9.456'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.717'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -3.92%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '62.02'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -12.99%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '5.149'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -13.57%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '6.857'
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -2.96%  '
